$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 65 (shifts existing rows 65-74 down to 66-75)
$ws.Rows.Item(65).Insert()

# Populate the new row 65 with the weekly price entry
$ws.Cells.Item(65, 1).Value = 11
$ws.Cells.Item(65, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(65, 3).Value = "Bíobío"
$ws.Cells.Item(65, 4).Value = 44476
$ws.Cells.Item(65, 5).Value = 8
$ws.Cells.Item(65, 6).Value = 100112043
$ws.Cells.Item(65, 7).Value = "Pepino ensalada"
$ws.Cells.Item(65, 8).Value = "Sin especificar"
$ws.Cells.Item(65, 9).Value = "Primera"
$ws.Cells.Item(65, 10).Value = 100
$ws.Cells.Item(65, 11).Value = 16000
$ws.Cells.Item(65, 12).Value = 17000
$ws.Cells.Item(65, 13).Value = 16500
$ws.Cells.Item(65, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(65, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(65, 16).Value = 275
$ws.Cells.Item(65, 17).Value = 60
$ws.Cells.Item(65, 18).Value = "Hortaliza"
